# Generate Report for handoff
#
# A new source file (711514bb-fefd-4675-bc73-e008cd09d5ab.md) failed its
# handoff transform, so it gets its own "Handoff transform failed" row on
# every sheet (pushing the ".localization-config" row down by one). The
# original source file's handoff re-ran with a new content hash / new
# timestamps.

$wb = $excel.ActiveWorkbook

$uuid1 = "0508241f-403d-4f18-9247-f5fe9e374d5f"
$uuid2 = "711514bb-fefd-4675-bc73-e008cd09d5ab"
$hash1 = "46d937e9d7caef2a89bb212c306488fe8ed2ff9d"

$srcCommit = "e6a29b81199da55149c6edee6e9722352571b627"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$uuid1.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("A3").Value = "$uuid2.md"
$ws1.Range("B3").Value = "Handoff transform failed"
$ws1.Range("C3").Value = "Handoff transform failed"

$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$uuid1.md", "", "", "$uuid1.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$uuid2.md", "", "", "$uuid2.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$xlf1zh = "$uuid1.$hash1.zh-cn.xlf"

$ws2.Range("A2").Value = "$uuid1.md"
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = $xlf1zh
$ws2.Range("D2").Value = "2016-01-18 06:54:25"
$ws2.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "$uuid2.md"
$ws2.Range("B3").Value = "Handoff transform failed"
$ws2.Range("D3").Value = "0001-01-01 00:00:00"
$ws2.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Ignored"

$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$uuid1.md", "", "", "$uuid1.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f1b97aaf310dd1ed33f96b469cf0aba6c2f6fbe9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$xlf1zh", "", "", $xlf1zh)
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$uuid2.md", "", "", "$uuid2.md")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$xlf1de = "$uuid1.$hash1.de-de.xlf"

$ws3.Range("A2").Value = "$uuid1.md"
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = $xlf1de
$ws3.Range("D2").Value = "2016-01-18 06:54:35"
$ws3.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "$uuid2.md"
$ws3.Range("B3").Value = "Handoff transform failed"
$ws3.Range("D3").Value = "0001-01-01 00:00:00"
$ws3.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Ignored"

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$uuid1.md", "", "", "$uuid1.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/897dc8757d92a67c843d757676d8277eb7886158/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$xlf1de", "", "", $xlf1de)
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$uuid2.md", "", "", "$uuid2.md")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/.localization-config", "", "", ".localization-config")

Write-Host "Report regenerated for handoff"
